$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before row 39 (shifts existing rows 39+ down by 2)
$ws.Rows.Item(39).Resize(2).Insert()

# Row 39: new "Primera" record for Provincia de Curicó dated 44980
$ws.Cells.Item(39, 1).Value = 9
$ws.Cells.Item(39, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(39, 3).Value = "Metropolitana"
$ws.Cells.Item(39, 4).Value = 44980
$ws.Cells.Item(39, 5).Value = 13
$ws.Cells.Item(39, 6).Value = "Fruta"
$ws.Cells.Item(39, 7).Value = 100101
$ws.Cells.Item(39, 8).Value = "Berries"
$ws.Cells.Item(39, 9).Value = 100101004
$ws.Cells.Item(39, 10).Value = "Frambuesa"
$ws.Cells.Item(39, 11).Value = "Sin especificar"
$ws.Cells.Item(39, 12).Value = "Primera"
$ws.Cells.Item(39, 13).Value = 300
$ws.Cells.Item(39, 14).Value = 6000
$ws.Cells.Item(39, 15).Value = 6000
$ws.Cells.Item(39, 16).Value = 6000
$ws.Cells.Item(39, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(39, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(39, 19).Value = 3000
$ws.Cells.Item(39, 20).Value = 2

# Row 40: new "Segunda" record for Provincia de Curicó dated 44980
$ws.Cells.Item(40, 1).Value = 9
$ws.Cells.Item(40, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(40, 3).Value = "Metropolitana"
$ws.Cells.Item(40, 4).Value = 44980
$ws.Cells.Item(40, 5).Value = 13
$ws.Cells.Item(40, 6).Value = "Fruta"
$ws.Cells.Item(40, 7).Value = 100101
$ws.Cells.Item(40, 8).Value = "Berries"
$ws.Cells.Item(40, 9).Value = 100101004
$ws.Cells.Item(40, 10).Value = "Frambuesa"
$ws.Cells.Item(40, 11).Value = "Sin especificar"
$ws.Cells.Item(40, 12).Value = "Segunda"
$ws.Cells.Item(40, 13).Value = 250
$ws.Cells.Item(40, 14).Value = 5000
$ws.Cells.Item(40, 15).Value = 5000
$ws.Cells.Item(40, 16).Value = 5000
$ws.Cells.Item(40, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(40, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(40, 19).Value = 2500
$ws.Cells.Item(40, 20).Value = 2

# Ensure date cells use the same number format as other date cells in column D
$ws.Cells.Item(39, 4).NumberFormat = $ws.Cells.Item(41, 4).NumberFormat
$ws.Cells.Item(40, 4).NumberFormat = $ws.Cells.Item(41, 4).NumberFormat
